# Update the credentials in row 2 of Sheet1:
#   A2: "mahesh23456"  -> "naveen2747"
#   B2: "India@123"    -> "Tester@2747"
# (rows 3-4 keep their existing values/order)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "naveen2747"
$ws.Range("B2").Value = "Tester@2747"

# Move the active selection to E3, matching the saved sheet view.
$ws.Range("E3").Select()
